$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp string (A1)
$ws.Range("A1").Value2 = "Datos actualizados a 25 de Agosto de 2020 a las 21:33"

# Update country rows with refreshed case counts (and relabel rows whose rank changed)
$ws.Range("A4").Value2 = "Estados Unidos"
$ws.Range("B4").Value2 = 5933627
$ws.Range("C4").Value2 = 17997
$ws.Range("D4").Value2 = 3225350
$ws.Range("E4").Value2 = 2526512
$ws.Range("F4").Value2 = 0
$ws.Range("G4").Value2 = 651
$ws.Range("H4").Value2 = 181765

$ws.Range("A23").Value2 = "Alemania"
$ws.Range("B23").Value2 = 236983
$ws.Range("C23").Value2 = 866
$ws.Range("D23").Value2 = 209600
$ws.Range("E23").Value2 = 18044
$ws.Range("F23").Value2 = 0
$ws.Range("G23").Value2 = 3
$ws.Range("H23").Value2 = 9339

$ws.Range("A59").Value2 = "Argelia"
$ws.Range("B59").Value2 = 42228
$ws.Range("C59").Value2 = 370
$ws.Range("D59").Value2 = 29369
$ws.Range("E59").Value2 = 11403
$ws.Range("F59").Value2 = 0
$ws.Range("G59").Value2 = 10
$ws.Range("H59").Value2 = 1456

$ws.Range("A76").Value2 = "Camerun"
$ws.Range("B76").Value2 = 18973
$ws.Range("C76").Value2 = 211
$ws.Range("D76").Value2 = 16540
$ws.Range("E76").Value2 = 2023
$ws.Range("F76").Value2 = 0
$ws.Range("G76").Value2 = 2
$ws.Range("H76").Value2 = 410

$ws.Range("A120").Value2 = "Mozambique"
$ws.Range("B120").Value2 = 3508
$ws.Range("C120").Value2 = 68
$ws.Range("D120").Value2 = 1809
$ws.Range("E120").Value2 = 1678
$ws.Range("F120").Value2 = 0
$ws.Range("G120").Value2 = 0
$ws.Range("H120").Value2 = 21

$ws.Range("A121").Value2 = "Eslovaquia"
$ws.Range("B121").Value2 = 3452
$ws.Range("C121").Value2 = 28
$ws.Range("D121").Value2 = 2167
$ws.Range("E121").Value2 = 1252
$ws.Range("F121").Value2 = 0
$ws.Range("G121").Value2 = 0
$ws.Range("H121").Value2 = 33

$ws.Range("A163").Value2 = "Republica del Chad"
$ws.Range("B163").Value2 = 995
$ws.Range("C163").Value2 = 8
$ws.Range("D163").Value2 = 871
$ws.Range("E163").Value2 = 47
$ws.Range("F163").Value2 = 0
$ws.Range("G163").Value2 = 1
$ws.Range("H163").Value2 = 77

$ws.Range("A182").Value2 = "Eritrea"
$ws.Range("B182").Value2 = 315
$ws.Range("C182").Value2 = 9
$ws.Range("D182").Value2 = 276
$ws.Range("E182").Value2 = 39
$ws.Range("F182").Value2 = 0
$ws.Range("G182").Value2 = 0
$ws.Range("H182").Value2 = 0

$ws.Range("A193").Value2 = "Monaco"
$ws.Range("B193").Value2 = 121
$ws.Range("C193").Value2 = 6
$ws.Range("D193").Value2 = 85
$ws.Range("E193").Value2 = 35
$ws.Range("F193").Value2 = 0
$ws.Range("G193").Value2 = 0
$ws.Range("H193").Value2 = 1

$ws.Range("A202").Value2 = "Islas Virgenes Britanicas"
$ws.Range("B202").Value2 = 26
$ws.Range("C202").Value2 = 5
$ws.Range("D202").Value2 = 8
$ws.Range("E202").Value2 = 17
$ws.Range("F202").Value2 = 0
$ws.Range("G202").Value2 = 0
$ws.Range("H202").Value2 = 1

$ws.Range("A204").Value2 = "Santa Lucia"
$ws.Range("B204").Value2 = 26
$ws.Range("C204").Value2 = 0
$ws.Range("D204").Value2 = 25
$ws.Range("E204").Value2 = 1
$ws.Range("F204").Value2 = 0
$ws.Range("G204").Value2 = 0
$ws.Range("H204").Value2 = 0

$ws.Range("A205").Value2 = "Granada"
$ws.Range("B205").Value2 = 24
$ws.Range("C205").Value2 = 0
$ws.Range("D205").Value2 = 24
$ws.Range("E205").Value2 = 0
$ws.Range("F205").Value2 = 0
$ws.Range("G205").Value2 = 0
$ws.Range("H205").Value2 = 0

$ws.Range("A206").Value2 = "Nueva Caledonia"
$ws.Range("B206").Value2 = 23
$ws.Range("C206").Value2 = 0
$ws.Range("D206").Value2 = 23
$ws.Range("E206").Value2 = 0
$ws.Range("F206").Value2 = 0
$ws.Range("G206").Value2 = 0
$ws.Range("H206").Value2 = 0

$ws.Range("A207").Value2 = "Laos"
$ws.Range("B207").Value2 = 22
$ws.Range("C207").Value2 = 0
$ws.Range("D207").Value2 = 21
$ws.Range("E207").Value2 = 1
$ws.Range("F207").Value2 = 0
$ws.Range("G207").Value2 = 0
$ws.Range("H207").Value2 = 0
